$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure rows: insert part-total subtotal rows, and re-space sections ---
# Insert blank row at 6 (will hold F6 = SUM(F3:F5), the Motor section subtotal)
$ws.Rows("6:6").Insert()
# Insert blank row before the (now shifted) pulley header row, so the pulley section
# occupies rows 9-12 and has room below it for a subtotal row
$ws.Rows("8:8").Insert()
# Insert a new row for the pulley section subtotal (F13 = SUM(F9:F12))
$ws.Rows("13:13").Insert()
# Remove one of the now-duplicated blank spacer rows
$ws.Rows("15:15").Delete()

# Copy formatting (grey column-A fill) from a neighboring blank styled cell into the
# two freshly inserted rows that need it
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Fill in the new subtotal formulas ---
$ws.Range("F6").Formula = "=SUM(F3:F5)"
$ws.Range("F13").Formula = "=SUM(F9:F12)"

# --- Bearings section: replace the "Frakt" shipping note with a flag value, and
#     turn the old shipping line into a flat shipping-cost total ---
$ws.Range("G22").Value = 1
$ws.Range("G24").ClearContents()
$ws.Range("F24").Value = 1830

# --- Grand total now sums the three section subtotals instead of one big SUM ---
$ws.Range("F31").Formula = "=F6+F13+F24"

Write-Host "Edit complete"
